# Auto-generated edit script: updates computed leve-profit cells per diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 48.57
$ws.Range("I15").Value = 48.57
$ws.Range("K15").Value = 145.71
$ws.Range("M15").Value = 23.28999999999999
$ws.Range("H98").Value = 5377.3486
$ws.Range("I98").Value = 3155.2
$ws.Range("J98").Value = 7309.6523
$ws.Range("K98").Value = 3155.2
$ws.Range("L98").Value = 7309.6523
$ws.Range("M98").Value = -1657.2
$ws.Range("N98").Value = -10305.6523
$ws.Range("H112").Value = 1337.8959
$ws.Range("J112").Value = 1337.8959
$ws.Range("L112").Value = 4013.6877
$ws.Range("N112").Value = -6229.6877
$ws.Range("H122").Value = 5377.3486
$ws.Range("I122").Value = 3155.2
$ws.Range("J122").Value = 7309.6523
$ws.Range("K122").Value = 9465.599999999999
$ws.Range("L122").Value = 21928.9569
$ws.Range("M122").Value = -7015.599999999999
$ws.Range("N122").Value = -26828.9569
$ws.Range("H132").Value = 205770.48
$ws.Range("I132").Value = 69511.664
$ws.Range("J132").Value = 911838.9399999999
$ws.Range("K132").Value = 208534.992
$ws.Range("L132").Value = 2735516.82
$ws.Range("M132").Value = -206004.992
$ws.Range("N132").Value = -2740576.82
$ws.Range("H137").Value = 621535.4
$ws.Range("I137").Value = 1987816.2
$ws.Range("J137").Value = 2842.17
$ws.Range("K137").Value = 5963448.6
$ws.Range("L137").Value = 8526.51
$ws.Range("M137").Value = -5960898.6
$ws.Range("N137").Value = -13626.51

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4315.5
$ws.Range("I74").Value = 4367
$ws.Range("K74").Value = 4367
$ws.Range("M74").Value = -3493
$ws.Range("H77").Value = 4315.5
$ws.Range("I77").Value = 4367
$ws.Range("K77").Value = 21835
$ws.Range("M77").Value = -17467
$ws.Range("H122").Value = 3207.4
$ws.Range("I122").Value = 2843.6562
$ws.Range("J122").Value = 4662.375
$ws.Range("K122").Value = 8530.9686
$ws.Range("L122").Value = 13987.125
$ws.Range("M122").Value = -6080.9686
$ws.Range("N122").Value = -18887.125
$ws.Range("H132").Value = 1637.9318
$ws.Range("I132").Value = 939.17145
$ws.Range("J132").Value = 4355.3335
$ws.Range("K132").Value = 2817.51435
$ws.Range("L132").Value = 13066.0005
$ws.Range("M132").Value = -287.5143500000004
$ws.Range("N132").Value = -18126.0005
$ws.Range("H133").Value = 34260.5
$ws.Range("J133").Value = 34260.5
$ws.Range("L133").Value = 34260.5
$ws.Range("N133").Value = -39320.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 437.66666
$ws.Range("I80").Value = 510.2
$ws.Range("J80").Value = 401.4
$ws.Range("K80").Value = 510.2
$ws.Range("L80").Value = 401.4
$ws.Range("M80").Value = 487.8
$ws.Range("N80").Value = -2397.4
$ws.Range("H83").Value = 437.66666
$ws.Range("I83").Value = 510.2
$ws.Range("J83").Value = 401.4
$ws.Range("K83").Value = 2551
$ws.Range("L83").Value = 2007
$ws.Range("M83").Value = 2441
$ws.Range("N83").Value = -11991
$ws.Range("H134").Value = 4077.2285
$ws.Range("I134").Value = 1418.6111
$ws.Range("K134").Value = 4255.8333
$ws.Range("M134").Value = -1720.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2506.6726
$ws.Range("I31").Value = 1159.579
$ws.Range("J31").Value = 3217.639
$ws.Range("K31").Value = 1159.579
$ws.Range("L31").Value = 3217.639
$ws.Range("M31").Value = -864.579
$ws.Range("N31").Value = -3807.639
$ws.Range("H34").Value = 2506.6726
$ws.Range("I34").Value = 1159.579
$ws.Range("J34").Value = 3217.639
$ws.Range("K34").Value = 1159.579
$ws.Range("L34").Value = 3217.639
$ws.Range("M34").Value = -957.579
$ws.Range("N34").Value = -3621.639
$ws.Range("H58").Value = 2820.742
$ws.Range("I58").Value = 1571.963
$ws.Range("J58").Value = 11250
$ws.Range("K58").Value = 1571.963
$ws.Range("L58").Value = 11250
$ws.Range("M58").Value = -1368.963
$ws.Range("N58").Value = -11656
$ws.Range("H105").Value = 1889.1666
$ws.Range("I105").Value = 1645.909
$ws.Range("J105").Value = 2271.4285
$ws.Range("K105").Value = 1645.909
$ws.Range("L105").Value = 2271.4285
$ws.Range("M105").Value = 101.0909999999999
$ws.Range("N105").Value = -5765.4285
$ws.Range("H136").Value = 2820.742
$ws.Range("I136").Value = 1571.963
$ws.Range("J136").Value = 11250
$ws.Range("K136").Value = 4715.889
$ws.Range("L136").Value = 33750
$ws.Range("M136").Value = -2165.889
$ws.Range("N136").Value = -38850

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1411.1014
$ws.Range("I68").Value = 1076.1818
$ws.Range("J68").Value = 1567.8723
$ws.Range("K68").Value = 3228.5454
$ws.Range("L68").Value = 4703.6169
$ws.Range("M68").Value = -2417.5454
$ws.Range("N68").Value = -6325.6169
$ws.Range("H71").Value = 1411.1014
$ws.Range("I71").Value = 1076.1818
$ws.Range("J71").Value = 1567.8723
$ws.Range("K71").Value = 9685.636200000001
$ws.Range("L71").Value = 14110.8507
$ws.Range("M71").Value = -5629.636200000001
$ws.Range("N71").Value = -22222.8507

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3430.56
$ws.Range("I126").Value = 2849.5144
$ws.Range("J126").Value = 4786.3335
$ws.Range("K126").Value = 8548.5432
$ws.Range("L126").Value = 14359.0005
$ws.Range("M126").Value = -6078.5432
$ws.Range("N126").Value = -19299.0005
$ws.Range("H132").Value = 4972.0557
$ws.Range("I132").Value = 2800
$ws.Range("J132").Value = 5243.5625
$ws.Range("K132").Value = 8400
$ws.Range("L132").Value = 15730.6875
$ws.Range("M132").Value = -5870
$ws.Range("N132").Value = -20790.6875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2410.6365
$ws.Range("I100").Value = 2020
$ws.Range("J100").Value = 2633.8572
$ws.Range("K100").Value = 2020
$ws.Range("L100").Value = 2633.8572
$ws.Range("M100").Value = -1479
$ws.Range("N100").Value = -3715.8572

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 22777.777
$ws.Range("J64").Value = 22777.777
$ws.Range("L64").Value = 22777.777
$ws.Range("N64").Value = -23273.777
$ws.Range("H67").Value = 22777.777
$ws.Range("J67").Value = 22777.777
$ws.Range("L67").Value = 22777.777
$ws.Range("N67").Value = -24493.777
$ws.Range("H136").Value = 4153.9614
$ws.Range("I136").Value = 1839.6154
$ws.Range("J136").Value = 6468.3076
$ws.Range("K136").Value = 5518.8462
$ws.Range("L136").Value = 19404.9228
$ws.Range("M136").Value = -2968.8462
$ws.Range("N136").Value = -24504.9228
